$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 187.55
$ws.Range("I33").Value = 124.066666
$ws.Range("J33").Value = 378
$ws.Range("K33").Value = 124.066666
$ws.Range("L33").Value = 378
$ws.Range("M33").Value = 104.933334
$ws.Range("N33").Value = -836
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 25880
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 25880
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 25880
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -26454
$ws.Range("H80").Value = 42259.6
$ws.Range("J80").Value = 42259.6
$ws.Range("L80").Value = 42259.6
$ws.Range("N80").Value = -44255.6
$ws.Range("H83").Value = 42259.6
$ws.Range("J83").Value = 42259.6
$ws.Range("L83").Value = 126778.8
$ws.Range("N83").Value = -136762.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 40000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 40000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 40000
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -40340
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 34950
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 34950
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 34950
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -35410
$ws.Range("H41").Value = 17395
$ws.Range("J41").Value = 19894
$ws.Range("L41").Value = 19894
$ws.Range("N41").Value = -20750
$ws.Range("H50").Value = 9100.200000000001
$ws.Range("J50").Value = 9100.200000000001
$ws.Range("L50").Value = 9100.200000000001
$ws.Range("N50").Value = -10350.2
$ws.Range("H51").Value = 9260.333000000001
$ws.Range("J51").Value = 9260.333000000001
$ws.Range("L51").Value = 9260.333000000001
$ws.Range("N51").Value = -10732.333
$ws.Range("H60").Value = 30516
$ws.Range("J60").Value = 30516
$ws.Range("L60").Value = 30516
$ws.Range("N60").Value = -31538
$ws.Range("H61").Value = 9260.333000000001
$ws.Range("J61").Value = 9260.333000000001
$ws.Range("L61").Value = 9260.333000000001
$ws.Range("N61").Value = -9956.333000000001
$ws.Range("H62").Value = 3922.2222
$ws.Range("I62").Value = 3892.3076
$ws.Range("K62").Value = 3892.3076
$ws.Range("M62").Value = -3268.3076
$ws.Range("H65").Value = 3922.2222
$ws.Range("I65").Value = 3892.3076
$ws.Range("K65").Value = 19461.538
$ws.Range("M65").Value = -16341.538
$ws.Range("H86").Value = 58826350
$ws.Range("J86").Value = 2655.4285
$ws.Range("L86").Value = 2655.4285
$ws.Range("N86").Value = -4901.4285
$ws.Range("H89").Value = 58826350
$ws.Range("J89").Value = 2655.4285
$ws.Range("L89").Value = 13277.1425
$ws.Range("N89").Value = -24509.1425
$ws.Range("H109").Value = 13714.286
$ws.Range("J109").Value = 13714.286
$ws.Range("L109").Value = 13714.286
$ws.Range("N109").Value = -15794.286
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 1400
$ws.Range("J62").Value = 3000
$ws.Range("L62").Value = 9000
$ws.Range("N62").Value = -10372
$ws.Range("H65").Value = 1400
$ws.Range("J65").Value = 3000
$ws.Range("L65").Value = 27000
$ws.Range("N65").Value = -33864
$ws.Range("H131").Value = 891.63
$ws.Range("I131").Value = 559.0909
$ws.Range("J131").Value = 932.73035
$ws.Range("K131").Value = 1677.2727
$ws.Range("L131").Value = 2798.19105
$ws.Range("M131").Value = 3362.7273
$ws.Range("N131").Value = -12878.19105
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 5183.3335
$ws.Range("J19").Value = 5750
$ws.Range("L19").Value = 5750
$ws.Range("N19").Value = -6326
$ws.Range("H46").Value = 11910
$ws.Range("J46").Value = 15756
$ws.Range("L46").Value = 15756
$ws.Range("N46").Value = -16068
$ws.Range("H123").Value = 34413.2
$ws.Range("J123").Value = 34413.2
$ws.Range("L123").Value = 34413.2
$ws.Range("N123").Value = -39313.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1274.6666
$ws.Range("I22").Value = 1166.6666
$ws.Range("J22").Value = 1310.6666
$ws.Range("K22").Value = 1166.6666
$ws.Range("L22").Value = 1310.6666
$ws.Range("M22").Value = -871.6666
$ws.Range("N22").Value = -1900.6666
$ws.Range("H27").Value = 1274.6666
$ws.Range("I27").Value = 1166.6666
$ws.Range("J27").Value = 1310.6666
$ws.Range("K27").Value = 1166.6666
$ws.Range("L27").Value = 1310.6666
$ws.Range("M27").Value = -1059.6666
$ws.Range("N27").Value = -1524.6666
$ws.Range("H68").Value = 1639.5312
$ws.Range("I68").Value = 1691.2963
$ws.Range("J68").Value = 1360
$ws.Range("K68").Value = 1691.2963
$ws.Range("L68").Value = 1360
$ws.Range("M68").Value = -942.2963
$ws.Range("N68").Value = -2858
$ws.Range("H71").Value = 1639.5312
$ws.Range("I71").Value = 1691.2963
$ws.Range("J71").Value = 1360
$ws.Range("K71").Value = 8456.4815
$ws.Range("L71").Value = 6800
$ws.Range("M71").Value = -4712.4815
$ws.Range("N71").Value = -14288
$ws.Range("H127").Value = 26864.166
$ws.Range("J127").Value = 26864.166
$ws.Range("L127").Value = 26864.166
$ws.Range("N127").Value = -36784.166
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 2000
$ws.Range("I52").Value = 2000
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 2000
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -1774
$ws.Range("N52").ClearContents()
$ws.Range("H62").Value = 6616.25
$ws.Range("I62").Value = 3922.5
$ws.Range("J62").Value = 9310
$ws.Range("K62").Value = 3922.5
$ws.Range("L62").Value = 9310
$ws.Range("M62").Value = -3298.5
$ws.Range("N62").Value = -10558
$ws.Range("H65").Value = 6616.25
$ws.Range("I65").Value = 3922.5
$ws.Range("J65").Value = 9310
$ws.Range("K65").Value = 19612.5
$ws.Range("L65").Value = 46550
$ws.Range("M65").Value = -16492.5
$ws.Range("N65").Value = -52790
$ws.Range("H109").Value = 18738.5
$ws.Range("J109").Value = 18738.5
$ws.Range("L109").Value = 18738.5
$ws.Range("N109").Value = -21512.5
$ws.Range("H124").Value = 40705
$ws.Range("J124").Value = 40705
$ws.Range("L124").Value = 40705
$ws.Range("N124").Value = -50525
